# The authored change swapped the contents of ppt/theme/theme1.xml and
# ppt/theme/theme2.xml: the "Office Theme" colour scheme and the
# "Integral" colour scheme traded places between the two theme parts
# (fonts/format scheme were already identical between the two themes).
#
# theme2.xml is the theme actually in force for the deck (it is the
# target of the slide master's theme relationship, so it drives
# Slide.ThemeColorScheme / Master.Theme for every slide). This script
# rewrites its 12 theme colours from the "Integral" palette to the
# "Office Theme" palette it is supposed to have after the edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$colorScheme = $s.ThemeColorScheme

function ToRgbLong([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Index order matches the OOXML <a:clrScheme> child order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = ToRgbLong $officeThemeColors[$i - 1]
}
